$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 117-120 (match ids 115-118 in column A) have their B and E:AD
# values cyclically rotated: new row117 <- old row118, new row118 <- old row119,
# new row119 <- old row120, new row120 <- old row117 (wrap-around).
# Columns B, E..AD map to column indices 2, 5..30.

$rows = @(117, 118, 119, 120)
$cols = @(2) + @(5..30)

# Capture the original values before overwriting anything.
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# new row[i] = old row[i+1], wrapping the last back to the first
for ($i = 0; $i -lt $rows.Length; $i++) {
    $destRow = $rows[$i]
    $srcRow = $rows[($i + 1) % $rows.Length]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $orig[$srcRow][$c]
    }
}
